# Add a "Site" column (E) to the Past_Sample_Data sheet that classifies
# each sample row as Inshore (IPa- rows, 2:111) or Offshore (OPa- rows, 112:221),
# then leave the selection on the newly-added Offshore block (to match the
# author's last on-screen selection after filling the column down).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("E1").Value = "Site"

# Inshore samples (rows 2-111, the "IPa-4-temp" experiment)
$ws.Range("E2:E111").Value = "Inshore"

# Offshore samples (rows 112-221, the "OPa-4-temp" experiment)
$ws.Range("E112:E221").Value = "Offshore"

# Match the author's final selection (the Offshore fill-down range)
$ws.Range("E112:E221").Select()
